$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the authoritative diff: (row, col, new value)
# Values are prefixed with a leading apostrophe marker so Excel stores them
# as literal text (matching the original inlineStr cells) instead of
# auto-converting number-like strings (e.g. "56.227.14", "1.00") into numbers.
$updates = @(
    @(2, 4, '56.227.14'),
    @(2, 5, '  +4.56%  '),
    @(3, 4, '2.474.27'),
    @(3, 5, '  +2.19%  '),
    @(4, 5, '  +0.07%  '),
    @(5, 4, '485.84'),
    @(5, 5, '  +5.24%  '),
    @(6, 4, '146.19'),
    @(6, 5, '  +11.47%  '),
    @(7, 4, '0.997'),
    @(7, 5, '  -0.23%  '),
    @(8, 5, '  +5.26%  '),
    @(9, 4, '2.484.63'),
    @(9, 5, '  +2.75%  '),
    @(10, 5, '  +8.23%  '),
    @(11, 4, '0.0962'),
    @(11, 5, '  +1.64%  '),
    @(12, 4, '0.329'),
    @(12, 5, '  +5.28%  '),
    @(13, 5, '  +1.55%  '),
    @(14, 4, '2.901.77'),
    @(14, 5, '  +2.18%  '),
    @(15, 4, '56.248.84'),
    @(15, 5, '  +4.55%  '),
    @(16, 4, '21.02'),
    @(16, 5, '  +7.24%  '),
    @(17, 5, '  +2.00%  '),
    @(18, 4, '2.485.58'),
    @(18, 5, '  +3.04%  '),
    @(19, 4, '4.49'),
    @(19, 5, '  +7.85%  '),
    @(20, 5, '  +6.03%  '),
    @(21, 4, '316.41'),
    @(21, 5, '  +3.29%  '),
    @(22, 5, '  -0.21%  '),
    @(23, 4, '5.77'),
    @(23, 5, '  +8.19%  '),
    @(24, 4, '58.29'),
    @(24, 5, '  +4.42%  '),
    @(25, 5, '  +6.45%  '),
    @(26, 4, '0.998'),
    @(26, 5, '  -1.01%  '),
    @(27, 5, '  +5.66%  '),
    @(28, 4, '2.587.70'),
    @(28, 5, '  +3.69%  '),
    @(29, 5, '  +6.98%  '),
    @(30, 4, '0.0₃0788'),
    @(30, 5, '  +10.13%  '),
    @(31, 4, '0.998'),
    @(31, 5, '  -0.13%  '),
    @(32, 4, '149.26'),
    @(32, 5, '  +2.42%  '),
    @(33, 5, '  +2.70%  '),
    @(34, 5, '  +5.44%  '),
    @(35, 5, '  +4.22%  '),
    @(36, 5, '  +8.08%  '),
    @(37, 5, '  +5.66%  '),
    @(38, 4, '0.857'),
    @(38, 5, '  +7.08%  '),
    @(39, 4, '34.12'),
    @(39, 5, '  +4.22%  '),
    @(40, 4, '3.49'),
    @(40, 5, '  +7.64%  '),
    @(41, 2, 'FirstDigitalUSD'),
    @(41, 3, 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'),
    @(41, 4, '0.995'),
    @(41, 5, '  -0.26%  '),
    @(42, 2, 'Hedera'),
    @(42, 3, 'https://coinranking.com/coin/jad286TjB+hedera-hbar'),
    @(42, 4, '0.0555'),
    @(42, 5, '  +6.12%  '),
    @(43, 4, '0.609'),
    @(43, 5, '  +2.76%  '),
    @(44, 5, '  +7.06%  '),
    @(45, 4, '4.74'),
    @(45, 5, '  +13.51%  '),
    @(46, 4, '0.0923'),
    @(46, 5, '  +6.55%  '),
    @(47, 4, '258.33'),
    @(47, 5, '  +15.53%  '),
    @(48, 5, '  +0.79%  '),
    @(49, 5, '  +5.21%  '),
    @(50, 2, 'EnergySwap'),
    @(50, 3, 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'),
    @(50, 4, '17.49'),
    @(50, 5, '  +6.52%  '),
    @(51, 2, 'Maker'),
    @(51, 3, 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'),
    @(51, 4, '1.877.94'),
    @(51, 5, '  -2.76%  '),
)

foreach ($u in $updates) {
    $r = $u[0]
    $c = $u[1]
    $v = $u[2]
    $ws.Cells.Item($r, $c).Value = "'" + $v
}

Write-Host "Applied $($updates.Count) cell updates."
